$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Words": append a new row (row 6) for the word "ਬਿਰਾਗੈਗੀ"
# ---------------------------------------------------------------------------
$wsWords = $wb.Worksheets.Item("Words")

$wsWords.Range("A6").Value = "ਬਿਰਾਗੈਗੀ"
$wsWords.Range("B6").Value = "ਬਿਰਾਗੈਗੀ"
$wsWords.Range("C6").Value = $true
$wsWords.Range("D6").Value = 45923.32848565972
$wsWords.Range("D6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWords.Range("E6").Value = $true
$wsWords.Range("F6").Value = 45923.32848565972
$wsWords.Range("F6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWords.Range("G6").Value = $true
$wsWords.Range("H6").Value = 45923.32852877315
$wsWords.Range("H6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWords.Range("I6").Value = $true
$wsWords.Range("J6").Value = 45923.3309484375
$wsWords.Range("J6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWords.Range("K6").Value = 0
$wsWords.Range("L6").Value = ""

# ---------------------------------------------------------------------------
# Sheet "Progress": mark rows 9, 10 and 79 as completed
# ---------------------------------------------------------------------------
$wsProgress = $wb.Worksheets.Item("Progress")

$wsProgress.Range("H9").Value = "completed"
$wsProgress.Range("I9").Value = 45923.33380170139
$wsProgress.Range("I9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsProgress.Range("J9").Value = 0

$wsProgress.Range("H10").Value = "completed"
$wsProgress.Range("I10").Value = 45923.33809645833
$wsProgress.Range("I10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsProgress.Range("J10").Value = 0

$wsProgress.Range("H79").Value = "completed"
$wsProgress.Range("I79").Value = 45923.33922759259
$wsProgress.Range("I79").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsProgress.Range("J79").Value = 0

# ---------------------------------------------------------------------------
# Sheet "Progress": append a new row (row 114) for the verse
# "ਮਨਿ ਬਿਰਾਗੈਗੀ ॥" of the word "ਬਿਰਾਗੈਗੀ"
# ---------------------------------------------------------------------------
$wsProgress.Range("A114").Value = "ਬਿਰਾਗੈਗੀ"
$wsProgress.Range("B114").Value = "ਬਿਰਾਗੈਗੀ"
$wsProgress.Range("C114").Value = ""
$wsProgress.Range("D114").Value = "ਮਨਿ ਬਿਰਾਗੈਗੀ ॥"
$wsProgress.Range("E114").Value = 1230
$wsProgress.Range("F114").Value = $true
$wsProgress.Range("G114").Value = 45923.32848565972
$wsProgress.Range("G114").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsProgress.Range("H114").Value = "completed"
$wsProgress.Range("I114").Value = 45923.33094693287
$wsProgress.Range("I114").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsProgress.Range("J114").Value = 0
$wsProgress.Range("K114").Value = ""
$wsProgress.Range("L114").Value = 45923.32848565972
$wsProgress.Range("L114").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsProgress.Range("M114").Value = ""
$wsProgress.Range("N114").Value = ""
$wsProgress.Range("O114").Value = ""
$wsProgress.Range("P114").Value = "ਮਨਿ ਬਿਰਾਗੈਗੀ"

Write-Output "edit applied"
